# Apply the commit's changes to the kobo "section2" workbook:
#  - settings!A2 (form_title) changes from the old duplicated title string
#    to "Any Regrets?"
#  - the "settings" sheet becomes the active/selected sheet (tab + view),
#    replacing "survey", and its selection moves to A2
#  - the "survey" sheet is no longer the tab-selected sheet

$wb = $excel.ActiveWorkbook

# --- settings sheet: update the form title text -----------------------
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A2").Value = "Any Regrets?"

# --- sheet/view selection: make "settings" the active tab -------------
# Activating the sheet and selecting A2 on it mirrors the diff's
# tabSelected/activeTab + <selection activeCell="A2"/> change, and also
# clears tabSelected from whichever sheet was previously active
# ("survey").
$settings.Activate()
$settings.Range("A2").Select()
